# "add index for list tag"
#
# The [row:list ...] loop tag gains an ", index" clause so the template
# engine exposes a loop index variable, and the column that used to print
# ${data.id} is switched to print that new ${index} value instead.
#
# Row 6: "[row:list datalist as data]"          -> "[row:list datalist as data, index]"
# Row 7: "${data.id}"                            -> "${index}"
# Row 8: "[/row:list]"                           (unchanged)
# Row 9: "${totalLabel}"                         (unchanged text)
#
# The active selection also moved from A6:B6 to A7:B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = '[row:list datalist as data, index]'
$ws.Range("A7").Value = '${index}'

$ws.Range("A7:B7").Select() | Out-Null
